# Fruta / hortaliza, semanal
# A new weekly price record is inserted at row 4 (pushing the existing
# rows 4-16 down to rows 5-17), and the new row is populated with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 4 - this shifts old rows 4..16 down to 5..17
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's data
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45050
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107011
$ws.Range("J4").Value = "Tuna"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("Q4").Value = "$/caja 18 kilos"
$ws.Range("R4").Value = "Región Metropolitana"
$ws.Range("S4").Value = 778
$ws.Range("T4").Value = 18
